$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Duplicate row 17 (values + full formatting, incl. the thick bottom
#    border) down onto the new row 18 -- this is currently the last row of
#    the table, and the new row becomes the new last row.
# ---------------------------------------------------------------------------
$ws.Range("A17:H17").Copy() | Out-Null
$ws.Range("A18:H18").PasteSpecial(-4104) | Out-Null          # xlPasteAll
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Write the new route's data into row 18.
# ---------------------------------------------------------------------------
$ws.Range("A18").Value = "POST"
$ws.Range("B18").Value = "/users/me/funds"
$ws.Range("C18").Value = "`n{amount, cardNum, expiryMonth, expiryYear, cvv, firstName, lastName, email, billingStreet, billingCity, billingState, billingCountryCode, billingZipCode}"
$ws.Range("D18").Value = "Authorization"
$ws.Range("E18").Value = "200, boolean"
$ws.Range("F18").Value = "200, boolean`n403"
$ws.Range("G18").ClearContents() | Out-Null
$ws.Range("H18").Value = "Adds funds to the current user by credit card"

# ---------------------------------------------------------------------------
# 3. Row 18 should use the lighter alternating-row shade (the same one used
#    by row 16) instead of the darker shade it inherited from row 17, while
#    keeping row 17's border layout (thick outer/bottom border). Columns E
#    and F keep their green/red success/error fills untouched.
# ---------------------------------------------------------------------------
$thin = 2
$medium = -4138
$lightShade = $ws.Range("A16").Interior.Color

function Set-RowBottomBorder($cell, $leftWeight, $rightWeight) {
    $cell.Borders.Item(7).LineStyle = 1
    $cell.Borders.Item(7).Weight = $leftWeight
    $cell.Borders.Item(8).LineStyle = 1
    $cell.Borders.Item(8).Weight = $thin
    $cell.Borders.Item(9).LineStyle = 1
    $cell.Borders.Item(9).Weight = $medium
    $cell.Borders.Item(10).LineStyle = 1
    $cell.Borders.Item(10).Weight = $rightWeight
}

foreach ($col in 1, 2, 3, 4, 7, 8) {
    $cell = $ws.Cells.Item(18, $col)
    $cell.Interior.Color = $lightShade
    $left = if ($col -eq 1) { $medium } else { $thin }
    $right = if ($col -eq 8) { $medium } else { $thin }
    Set-RowBottomBorder $cell $left $right
}

# ---------------------------------------------------------------------------
# 4. Row height -- tall enough to show the wrapped request-body text.
# ---------------------------------------------------------------------------
$ws.Rows(18).RowHeight = 124.5

# ---------------------------------------------------------------------------
# 5. Extend the autofilter / table range from A1:H17 to A1:H18.
# ---------------------------------------------------------------------------
$ws.AutoFilterMode = $false
$ws.Range("A1:H18").AutoFilter() | Out-Null

# ---------------------------------------------------------------------------
# 6. Keep the _FilterDatabase defined name lined up with the autofilter.
# ---------------------------------------------------------------------------
foreach ($i in 1..$wb.Names.Count) {
    $n = $wb.Names.Item($i)
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Routes!`$A`$1:`$H`$18"
    }
}

# ---------------------------------------------------------------------------
# 7. Move the active selection down to the newly added row, like a user
#    would after typing the new data in.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A18").Select()
